$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 85, shifting rows 85-99 down to 86-100.
$ws.Rows.Item(85).Insert()

# Set the new row's B85 cell to the new label.
$ws.Range("B85").Value = "Multiple Choice Questions"

# Update view state (selection) to match the post-edit workbook state.
[void]$ws.Range("F112").Select()
